$d = $word.ActiveDocument

# Paragraph: "Bonus 2" - top-level bullet (ilvl 0) in the same list (numId 3)
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Bonus 2"
$p1.Range.ListFormat.ListLevelNumber = 1

# Paragraph: mean/median explanation - sub-level bullet (ilvl 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "The mean summarizes the data in a more meaningful way. The mean takes into account the number of backers which can greatly vary from project to project depending on the Goal Amount.  The median is only the middle point which can vary greatly depending on the campaign."
$p2.Range.ListFormat.ListLevelNumber = 2

# Paragraph: variability explanation - sub-level bullet (ilvl 1)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "There is greater variability in successful campaigns. This makes sense as number of backers and Goal Amount can greatly vary in a successful campaign. A single person could pledge `$50 to make a campaign successful and at the same time 10,000 people could pledge &50 to make another campaign successful."
$p3.Range.ListFormat.ListLevelNumber = 2
